$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.178962333333333
$ws.Range("H2").Value = 3.536887
$ws.Range("I2").Value = 0.001182125215344215
$ws.Range("J2").Value = 0.001182125215344214
$ws.Range("M2").Value = 2.582049666666667
$ws.Range("N2").Value = 7.746149
$ws.Range("O2").Value = 0.02216861605835153
$ws.Range("P2").Value = 0.02216861605835152
$ws.Range("Q2").Value = 3.044139299795889
$ws.Range("R2").Value = 27.397253698163
$ws.Range("S2").Value = 0.00002620608003186201
$ws.Range("T2").Value = 0.000026206080031862
$ws.Range("G3").Value = 1.178962333333333
$ws.Range("H3").Value = 3.536887
$ws.Range("I3").Value = 0.001182125215344215
$ws.Range("J3").Value = 0.001182125215344214
$ws.Range("M3").Value = 93.97803500000002
$ws.Range("N3").Value = 281.934105
$ws.Range("O3").Value = 0.8068640207540504
$ws.Range("P3").Value = 0.8068640207540503
$ws.Range("Q3").Value = 110.7965634256817
$ws.Range("R3").Value = 997.1690708311352
$ws.Range("S3").Value = 0.0009538143042873808
$ws.Range("T3").Value = 0.0009538143042873802
$ws.Range("G4").Value = 1.178962333333333
$ws.Range("H4").Value = 3.536887
$ws.Range("I4").Value = 0.001182125215344215
$ws.Range("J4").Value = 0.001182125215344214
$ws.Range("M4").Value = 19.913116
$ws.Range("N4").Value = 59.73934799999999
$ws.Range("O4").Value = 0.1709673631875981
$ws.Range("P4").Value = 0.1709673631875981
$ws.Range("Q4").Value = 23.47681370329733
$ws.Range("R4").Value = 211.291323329676
$ws.Range("S4").Value = 0.000202104831024972
$ws.Range("T4").Value = 0.0002021048310249719
$ws.Range("G5").Value = 950.6576336666667
$ws.Range("H5").Value = 2851.972901
$ws.Range("I5").Value = 0.9532080272144655
$ws.Range("J5").Value = 0.9532080272144653
$ws.Range("M5").Value = 2.582049666666667
$ws.Range("N5").Value = 7.746149
$ws.Range("O5").Value = 0.02216861605835153
$ws.Range("P5").Value = 0.02216861605835152
$ws.Range("Q5").Value = 2454.645226123139
$ws.Range("R5").Value = 22091.80703510825
$ws.Range("S5").Value = 0.02113130277905618
$ws.Range("T5").Value = 0.02113130277905617
$ws.Range("G6").Value = 950.6576336666667
$ws.Range("H6").Value = 2851.972901
$ws.Range("I6").Value = 0.9532080272144655
$ws.Range("J6").Value = 0.9532080272144653
$ws.Range("M6").Value = 93.97803500000002
$ws.Range("N6").Value = 281.934105
$ws.Range("O6").Value = 0.8068640207540504
$ws.Range("P6").Value = 0.8068640207540503
$ws.Range("Q6").Value = 89340.9363697432
$ws.Range("R6").Value = 804068.4273276888
$ws.Range("S6").Value = 0.7691092614532999
$ws.Range("T6").Value = 0.7691092614532996
$ws.Range("G7").Value = 950.6576336666667
$ws.Range("H7").Value = 2851.972901
$ws.Range("I7").Value = 0.9532080272144655
$ws.Range("J7").Value = 0.9532080272144653
$ws.Range("M7").Value = 19.913116
$ws.Range("N7").Value = 59.73934799999999
$ws.Range("O7").Value = 0.1709673631875981
$ws.Range("P7").Value = 0.1709673631875981
$ws.Range("Q7").Value = 18930.55573548984
$ws.Range("R7").Value = 170375.0016194085
$ws.Range("S7").Value = 0.1629674629821094
$ws.Range("T7").Value = 0.1629674629821094
$ws.Range("G8").Value = 45.48781433333333
$ws.Range("H8").Value = 136.463443
$ws.Range("I8").Value = 0.04560984757019037
$ws.Range("J8").Value = 0.04560984757019036
$ws.Range("M8").Value = 2.582049666666667
$ws.Range("N8").Value = 7.746149
$ws.Range("O8").Value = 0.02216861605835153
$ws.Range("P8").Value = 0.02216861605835152
$ws.Range("Q8").Value = 117.4517958367785
$ws.Range("R8").Value = 1057.066162531007
$ws.Range("S8").Value = 0.001011107199263488
$ws.Range("T8").Value = 0.001011107199263487
$ws.Range("G9").Value = 45.48781433333333
$ws.Range("H9").Value = 136.463443
$ws.Range("I9").Value = 0.04560984757019037
$ws.Range("J9").Value = 0.04560984757019036
$ws.Range("M9").Value = 93.97803500000002
$ws.Range("N9").Value = 281.934105
$ws.Range("O9").Value = 0.8068640207540504
$ws.Range("P9").Value = 0.8068640207540503
$ws.Range("Q9").Value = 4274.855407491502
$ws.Range("R9").Value = 38473.69866742352
$ws.Range("S9").Value = 0.03680094499646316
$ws.Range("T9").Value = 0.03680094499646314
$ws.Range("G10").Value = 45.48781433333333
$ws.Range("H10").Value = 136.463443
$ws.Range("I10").Value = 0.04560984757019037
$ws.Range("J10").Value = 0.04560984757019036
$ws.Range("M10").Value = 19.913116
$ws.Range("N10").Value = 59.73934799999999
$ws.Range("O10").Value = 0.1709673631875981
$ws.Range("P10").Value = 0.1709673631875981
$ws.Range("Q10").Value = 905.8041234061291
$ws.Range("R10").Value = 8152.237110655162
$ws.Range("S10").Value = 0.007797795374463726
$ws.Range("T10").Value = 0.007797795374463725
